$wb = $excel.ActiveWorkbook

# Update the price for the first article on the "Lista" sheet
$listaSheet = $wb.Worksheets.Item("Lista")
$listaSheet.Range("C1").Value = 2003

# Make "Lista" the active sheet/tab and select cell C1 on it
$listaSheet.Activate()
$listaSheet.Range("C1").Select()
